$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.316967666666667
$ws.Range("H2").Value = 3.950903
$ws.Range("I2").Value = 0.02034337776957547
$ws.Range("J2").Value = 0.02034337776957546
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 2.442429333333334
$ws.Range("N2").Value = 7.327288
$ws.Range("O2").Value = 0.08913295894744963
$ws.Range("P2").Value = 0.08913295894744963
$ws.Range("Q2").Value = 3.216600460118222
$ws.Range("R2").Value = 28.949404141064
$ws.Range("S2").Value = 0.001813265455588029
$ws.Range("T2").Value = 0.001813265455588029

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.316967666666667
$ws.Range("H3").Value = 3.950903
$ws.Range("I3").Value = 0.02034337776957547
$ws.Range("J3").Value = 0.02034337776957546
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 15.82990933333333
$ws.Range("N3").Value = 47.489728
$ws.Range("O3").Value = 0.5776898596383203
$ws.Range("P3").Value = 0.5776898596383203
$ws.Range("Q3").Value = 20.84747875826489
$ws.Range("R3").Value = 187.627308824384
$ws.Range("S3").Value = 0.01175216304827538
$ws.Range("T3").Value = 0.01175216304827537

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.316967666666667
$ws.Range("H4").Value = 3.950903
$ws.Range("I4").Value = 0.02034337776957547
$ws.Range("J4").Value = 0.02034337776957546
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 9.129750999999999
$ws.Range("N4").Value = 27.389253
$ws.Range("O4").Value = 0.3331771814142301
$ws.Range("P4").Value = 0.3331771814142301
$ws.Range("Q4").Value = 12.02358687171766
$ws.Range("R4").Value = 108.212281845459
$ws.Range("S4").Value = 0.00677794926571206
$ws.Range("T4").Value = 0.006777949265712058

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 43.63696533333334
$ws.Range("H5").Value = 130.910896
$ws.Range("I5").Value = 0.6740661088064187
$ws.Range("J5").Value = 0.6740661088064185
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 2.442429333333334
$ws.Range("N5").Value = 7.327288
$ws.Range("O5").Value = 0.08913295894744963
$ws.Range("P5").Value = 0.08913295894744963
$ws.Range("Q5").Value = 106.5802041477831
$ws.Range("R5").Value = 959.2218373300481
$ws.Range("S5").Value = 0.06008150680410963
$ws.Range("T5").Value = 0.06008150680410961

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 43.63696533333334
$ws.Range("H6").Value = 130.910896
$ws.Range("I6").Value = 0.6740661088064187
$ws.Range("J6").Value = 0.6740661088064185
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 15.82990933333333
$ws.Range("N6").Value = 47.489728
$ws.Range("O6").Value = 0.5776898596383203
$ws.Range("P6").Value = 0.5776898596383203
$ws.Range("Q6").Value = 690.7692048084765
$ws.Range("R6").Value = 6216.922843276288
$ws.Range("S6").Value = 0.3894011557833288
$ws.Range("T6").Value = 0.3894011557833286

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 43.63696533333334
$ws.Range("H7").Value = 130.910896
$ws.Range("I7").Value = 0.6740661088064187
$ws.Range("J7").Value = 0.6740661088064185
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 9.129750999999999
$ws.Range("N7").Value = 27.389253
$ws.Range("O7").Value = 0.3331771814142301
$ws.Range("P7").Value = 0.3331771814142301
$ws.Range("Q7").Value = 398.3946278889653
$ws.Range("R7").Value = 3585.551651000688
$ws.Range("S7").Value = 0.2245834462189803
$ws.Range("T7").Value = 0.2245834462189802

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 19.78298933333333
$ws.Range("H8").Value = 59.348968
$ws.Range("I8").Value = 0.305590513424006
$ws.Range("J8").Value = 0.3055905134240059
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 2.442429333333334
$ws.Range("N8").Value = 7.327288
$ws.Range("O8").Value = 0.08913295894744963
$ws.Range("P8").Value = 0.08913295894744963
$ws.Range("Q8").Value = 48.31855344875378
$ws.Range("R8").Value = 434.866981038784
$ws.Range("S8").Value = 0.02723818668775198
$ws.Range("T8").Value = 0.02723818668775198

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 19.78298933333333
$ws.Range("H9").Value = 59.348968
$ws.Range("I9").Value = 0.305590513424006
$ws.Range("J9").Value = 0.3055905134240059
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 15.82990933333333
$ws.Range("N9").Value = 47.489728
$ws.Range("O9").Value = 0.5776898596383203
$ws.Range("P9").Value = 0.5776898596383203
$ws.Range("Q9").Value = 313.1629274889671
$ws.Range("R9").Value = 2818.466347400704
$ws.Range("S9").Value = 0.1765365408067163
$ws.Range("T9").Value = 0.1765365408067162

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 19.78298933333333
$ws.Range("H10").Value = 59.348968
$ws.Range("I10").Value = 0.305590513424006
$ws.Range("J10").Value = 0.3055905134240059
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 9.129750999999999
$ws.Range("N10").Value = 27.389253
$ws.Range("O10").Value = 0.3331771814142301
$ws.Range("P10").Value = 0.3331771814142301
$ws.Range("Q10").Value = 180.6137666489893
$ws.Range("R10").Value = 1625.523899840904
$ws.Range("S10").Value = 0.1018157859295378
$ws.Range("T10").Value = 0.1018157859295377
